$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 326, shifting existing rows 326:431 down to 327:432
$ws.Rows.Item(326).Insert()

# Populate the newly inserted row 326 with the new data record
$ws.Range("A326").Value = 7
$ws.Range("B326").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C326").Value = "Ñuble"
$ws.Range("D326").Value = 45015
$ws.Range("E326").Value = 16
$ws.Range("F326").Value = "Fruta"
$ws.Range("G326").Value = 100103
$ws.Range("H326").Value = "Frutos de hueso (carozo)"
$ws.Range("I326").Value = 100103004
$ws.Range("J326").Value = "Durazno"
$ws.Range("K326").Value = "Phillips Cling"
$ws.Range("L326").Value = "Primera"
$ws.Range("M326").Value = 70
$ws.Range("N326").Value = 16000
$ws.Range("O326").Value = 17000
$ws.Range("P326").Value = 16571
$ws.Range("Q326").Value = "`$/caja 16 kilos empedrada"
$ws.Range("R326").Value = "Región de O'Higgins"
$ws.Range("S326").Value = 1036
$ws.Range("T326").Value = 16
